$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in title cell A1
$ws.Range("A1").Value = "Datos actualizados a 16 de Agosto de 2020 a las 22:34"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 5558808
$ws.Range("C4").Value = 29019
$ws.Range("D4").Value = 2911515
$ws.Range("E4").Value = 2474324
$ws.Range("G4").Value = 363
$ws.Range("H4").Value = 172969

# Row 8: Sudafrica
$ws.Range("A8").Value = "Sudafrica"
$ws.Range("B8").Value = 587345
$ws.Range("C8").Value = 3692
$ws.Range("D8").Value = 472377
$ws.Range("E8").Value = 103129
$ws.Range("G8").Value = 162
$ws.Range("H8").Value = 11839

# Row 22: Alemania
$ws.Range("A22").Value = "Alemania"
$ws.Range("B22").Value = 224997
$ws.Range("C22").Value = 519
$ws.Range("D22").Value = 202900
$ws.Range("E22").Value = 12807
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 9290

# Row 68: Costa Rica
$ws.Range("A68").Value = "Costa Rica"
$ws.Range("B68").Value = 28465
$ws.Range("C68").Value = 728
$ws.Range("D68").Value = 9062
$ws.Range("E68").Value = 19109
$ws.Range("G68").Value = 3
$ws.Range("H68").Value = 294

# Row 102: Mauritania
$ws.Range("A102").Value = "Mauritania"
$ws.Range("B102").Value = 6701
$ws.Range("C102").Value = 8
$ws.Range("D102").Value = 5985
$ws.Range("E102").Value = 559
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 157

# Row 106: Zimbabue
$ws.Range("A106").Value = "Zimbabue"
$ws.Range("B106").Value = 5261
$ws.Range("C106").Value = 85
$ws.Range("D106").Value = 2092
$ws.Range("E106").Value = 3037
$ws.Range("G106").Value = 2
$ws.Range("H106").Value = 132

# Row 112: Namibia
$ws.Range("A112").Value = "Namibia"
$ws.Range("B112").Value = 4154
$ws.Range("C112").Value = 247
$ws.Range("D112").Value = 2370
$ws.Range("E112").Value = 1749
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 35

# Row 113: Nicaragua
$ws.Range("A113").Value = "Nicaragua"
$ws.Range("B113").Value = 4115
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 2913
$ws.Range("E113").Value = 1074
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 128

# Row 114: Montenegro
$ws.Range("A114").Value = "Montenegro"
$ws.Range("B114").Value = 3960
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 2830
$ws.Range("E114").Value = 1055
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 75

# Row 115: Suazilandia
$ws.Range("A115").Value = "Suazilandia"
$ws.Range("B115").Value = 3839
$ws.Range("C115").Value = 94
$ws.Range("D115").Value = 2268
$ws.Range("E115").Value = 1501
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 70

# Row 116: Congo
$ws.Range("A116").Value = "Congo"
$ws.Range("B116").Value = 3831
$ws.Range("C116").Value = 86
$ws.Range("D116").Value = 1625
$ws.Range("E116").Value = 2130
$ws.Range("G116").Value = 16
$ws.Range("H116").Value = 76

# Row 120: Cabo Verde
$ws.Range("A120").Value = "Cabo Verde"
$ws.Range("B120").Value = 3179
$ws.Range("C120").Value = 16
$ws.Range("D120").Value = 2317
$ws.Range("E120").Value = 827
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 35

# Row 126: Mali
$ws.Range("A126").Value = "Mali"
$ws.Range("B126").Value = 2640
$ws.Range("C126").Value = 26
$ws.Range("D126").Value = 1987
$ws.Range("E126").Value = 528
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 125

# Row 136: Sierra Leona
$ws.Range("A136").Value = "Sierra Leona"
$ws.Range("B136").Value = 1956
$ws.Range("C136").Value = 2
$ws.Range("D136").Value = 1506
$ws.Range("E136").Value = 381
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 69

# Row 138: Gambia
$ws.Range("A138").Value = "Gambia"
$ws.Range("B138").Value = 1872
$ws.Range("C138").Value = 183
$ws.Range("D138").Value = 401
$ws.Range("E138").Value = 1408
$ws.Range("G138").Value = 9
$ws.Range("H138").Value = 63

# Row 139: Yemen
$ws.Range("A139").Value = "Yemen"
$ws.Range("B139").Value = 1869
$ws.Range("C139").Value = 11
$ws.Range("D139").Value = 1013
$ws.Range("E139").Value = 326
$ws.Range("G139").Value = 2
$ws.Range("H139").Value = 530

# Row 142: Uganda
$ws.Range("A142").Value = "Uganda"
$ws.Range("B142").Value = 1500
$ws.Range("C142").Value = 66
$ws.Range("D142").Value = 1142
$ws.Range("E142").Value = 345
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 13

# Row 213: Islas Malvinas
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

# Row 214: Montserrat
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
